# daily auto push: 2026-02-16 14:10 UTC
# Two new readings were logged for 2026/02/16 (time slots 20 and 22),
# pushing every subsequent row down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right above the first "2026/12/29" row (row 822),
# shifting all existing rows (822-863) down to (824-865).
$ws.Range("A822:D823").EntireRow.Insert()

# Keep the date column as plain text (matches the rest of column A) instead
# of letting it auto-convert to a date serial number.
$ws.Range("A822:A823").NumberFormat = "@"

$ws.Range("A822").Value = "2026/02/16"
$ws.Range("B822").Value = "月"
$ws.Range("C822").Value = 20
$ws.Range("D822").Value = 201

$ws.Range("A823").Value = "2026/02/16"
$ws.Range("B823").Value = "月"
$ws.Range("C823").Value = 22
$ws.Range("D823").Value = 201
